# This script updates the "cryptos" price table (Sheet1, A1:E51) on
# $ws to reflect a refreshed data pull, as produced by the scheduled
# GitHub Actions job. A new coin (OKB) was inserted after "Polygon",
# shifting subsequent rows down by one and dropping "Aave" off the
# bottom of the fixed-size A1:E51 range; Price/Volume(1h) values were
# refreshed for every remaining row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' (price/volume refresh)
$ws.Range("D2").Value = '27.993.41'
$ws.Range("E2").Value = '  -0.54%  '

# Row 3: 'Ethereum' (price/volume refresh)
$ws.Range("D3").Value = '1.849.46'
$ws.Range("E3").Value = '  -1.23%  '

# Row 4: 'TetherUSD' (price/volume refresh)
$ws.Range("D4").Value = "'" + '0.9970'
$ws.Range("E4").Value = '  -1.00%  '

# Row 5: 'BNB' (price/volume refresh)
$ws.Range("D5").Value = "'" + '310.97'
$ws.Range("E5").Value = '  -0.85%  '

# Row 6: 'USDC' (price/volume refresh)
$ws.Range("D6").Value = "'" + '0.9978'
$ws.Range("E6").Value = '  -0.69%  '

# Row 7: 'XRP' (price/volume refresh)
$ws.Range("D7").Value = "'" + '0.5063'
$ws.Range("E7").Value = '  -1.52%  '

# Row 8: 'Cardano' (price/volume refresh)
$ws.Range("D8").Value = "'" + '0.3891'
$ws.Range("E8").Value = '  +0.26%  '

# Row 9: 'Dogecoin' (price/volume refresh)
$ws.Range("D9").Value = "'" + '0.08222'
$ws.Range("E9").Value = '  -1.92%  '

# Row 10: 'Polygon' (price/volume refresh)
$ws.Range("D10").Value = "'" + '1.104'
$ws.Range("E10").Value = '  -0.93%  '

# Row 11: 'Polkadot' -> 'OKB'
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = "'" + '41.32'
$ws.Range("E11").Value = '  -0.97%  '

# Row 12: 'WrappedEther' -> 'Polkadot'
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = "'" + '6.173'
$ws.Range("E12").Value = '  -0.38%  '

# Row 13: 'Solana' -> 'WrappedEther'
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.846.03'
$ws.Range("E13").Value = '  -1.72%  '

# Row 14: 'Chainlink' -> 'Solana'
$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").Value = "'" + '20.11'
$ws.Range("E14").Value = '  -2.28%  '

# Row 15: 'BinanceUSD' -> 'Chainlink'
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = "'" + '7.156'
$ws.Range("E15").Value = '  -1.95%  '

# Row 16: 'ShibaInu' -> 'BinanceUSD'
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = "'" + '0.9964'
$ws.Range("E16").Value = '  -1.08%  '

# Row 17: 'Litecoin' (price/volume refresh)
$ws.Range("D17").Value = "'" + '90.74'
$ws.Range("E17").Value = '  -0.30%  '

# Row 18: 'TRON' -> 'ShibaInu'
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = "'" + '0.00001090'
$ws.Range("E18").Value = '  -1.51%  '

# Row 19: 'Avalanche' -> 'TRON'
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = "'" + '0.06629'
$ws.Range("E19").Value = '  -0.68%  '

# Row 20: 'Dai' -> 'Avalanche'
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = "'" + '17.48'
$ws.Range("E20").Value = '  -1.26%  '

# Row 21: 'Uniswap' -> 'Dai'
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = "'" + '0.9990'
$ws.Range("E21").Value = '  -0.49%  '

# Row 22: 'WrappedBTC' -> 'Uniswap'
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = "'" + '5.896'
$ws.Range("E22").Value = '  -2.31%  '

# Row 23: 'Cosmos' -> 'WrappedBTC'
$ws.Range("B23").Value = 'WrappedBTC'
$ws.Range("C23").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D23").Value = '28.033.50'
$ws.Range("E23").Value = '  -0.50%  '

# Row 24: 'Toncoin' -> 'Cosmos'
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = "'" + '11.00'
$ws.Range("E24").Value = '  -1.00%  '

# Row 25: 'WrappedliquidstakedEther2.0' -> 'Toncoin'
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = "'" + '2.216'
$ws.Range("E25").Value = '  -1.41%  '

# Row 26: 'Monero' -> 'WrappedliquidstakedEther2.0'
$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.065.21'
$ws.Range("E26").Value = '  -0.92%  '

# Row 27: 'EthereumClassic' -> 'Monero'
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = "'" + '158.65'
$ws.Range("E27").Value = '  +0.03%  '

# Row 28: 'LidoDAOToken' -> 'EthereumClassic'
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = "'" + '20.46'
$ws.Range("E28").Value = '  -0.60%  '

# Row 29: 'BitcoinCash' -> 'LidoDAOToken'
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = "'" + '2.394'
$ws.Range("E29").Value = '  -3.28%  '

# Row 30: 'Stellar' -> 'BitcoinCash'
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = "'" + '125.64'
$ws.Range("E30").Value = '  +0.57%  '

# Row 31: 'ImmutableX' -> 'Stellar'
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = "'" + '0.1047'
$ws.Range("E31").Value = '  -1.43%  '

# Row 32: 'Filecoin' -> 'ImmutableX'
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = "'" + '1.026'
$ws.Range("E32").Value = '  -1.26%  '

# Row 33: 'HuobiToken' -> 'Filecoin'
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = "'" + '5.784'
$ws.Range("E33").Value = '  -1.72%  '

# Row 34: 'VeChain' -> 'HuobiToken'
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = "'" + '3.568'
$ws.Range("E34").Value = '  -0.67%  '

# Row 35: 'Hedera' -> 'VeChain'
$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").Value = "'" + '0.02416'
$ws.Range("E35").Value = '  -0.94%  '

# Row 36: 'FraxShare' -> 'Hedera'
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = "'" + '0.06419'
$ws.Range("E36").Value = '  -1.85%  '

# Row 37: 'Algorand' -> 'FraxShare'
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").Value = "'" + '9.000'
$ws.Range("E37").Value = '  -6.32%  '

# Row 38: 'TrustWalletToken' -> 'Algorand'
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = "'" + '0.2151'
$ws.Range("E38").Value = '  -1.61%  '

# Row 39: 'TheSandbox' -> 'TrustWalletToken'
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = "'" + '1.241'
$ws.Range("E39").Value = '  +0.80%  '

# Row 40: 'ARBITRUM' -> 'TheSandbox'
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = "'" + '0.6387'
$ws.Range("E40").Value = '  -1.69%  '

# Row 41: 'InternetComputer(DFINITY)' -> 'ARBITRUM'
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = "'" + '1.170'
$ws.Range("E41").Value = '  -2.98%  '

# Row 42: 'Aptos' -> 'InternetComputer(DFINITY)'
$ws.Range("B42").Value = 'InternetComputer(DFINITY)'
$ws.Range("C42").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D42").Value = "'" + '4.919'
$ws.Range("E42").Value = '  -1.76%  '

# Row 43: 'Decentraland' -> 'Aptos'
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = "'" + '10.99'
$ws.Range("E43").Value = '  -3.08%  '

# Row 44: 'EnergySwap' -> 'Decentraland'
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = "'" + '0.5970'
$ws.Range("E44").Value = '  -1.93%  '

# Row 45: 'PancakeSwap' -> 'EnergySwap'
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = "'" + '12.88'
$ws.Range("E45").Value = '  -1.07%  '

# Row 46: 'WEMIXTOKEN' -> 'PancakeSwap'
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = "'" + '3.639'
$ws.Range("E46").Value = '  -1.04%  '

# Row 47: 'NEARProtocol' -> 'WEMIXTOKEN'
$ws.Range("B47").Value = 'WEMIXTOKEN'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").Value = "'" + '1.259'
$ws.Range("E47").Value = '  -1.47%  '

# Row 48: 'EOS' -> 'NEARProtocol'
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = "'" + '1.987'
$ws.Range("E48").Value = '  -1.10%  '

# Row 49: 'Quant' -> 'EOS'
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = "'" + '1.195'
$ws.Range("E49").Value = '  -1.78%  '

# Row 50: 'Cronos' -> 'Quant'
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = "'" + '120.09'
$ws.Range("E50").Value = '  -1.13%  '

# Row 51: 'Aave' -> 'Cronos'
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'" + '0.06833'
$ws.Range("E51").Value = '  -0.63%  '
